$ErrorActionPreference = "Stop"
$d = $word.ActiveDocument

# Locate the paragraph "Str.strip([chars]) (Both beginning and ending whitespace)"
# which currently carries the _GoBack bookmark at its end.
$count = $d.Paragraphs.Count
$idx = -1
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "Str.strip(*beginning and ending whitespace)*") {
        $idx = $i
        break
    }
}
if ($idx -eq -1) {
    throw "Could not locate anchor paragraph (Str.strip ... whitespace))"
}

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Insert a fresh empty paragraph right after the anchor paragraph; we'll
# replace its contents with the full block of new paragraphs (ending in
# the paragraph that now owns the _GoBack bookmark).
$anchor = $d.Paragraphs.Item($idx)
$r = $anchor.Range
$r.Collapse(0)
$null = $r.InsertParagraphAfter()

$newBlockXml = @'
<w:body xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
    <w:p>
      <w:r>
        <w:t>Str.upper()</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Search and replace can be done through “replace” function, replace(</w:t>
      </w:r>
      <w:r>
        <w:t>Variable that needs to be replaced, what you want it replaced with)</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Chapter 7 Files:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>-A text file can be thought of as a sequence of lines</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">-Before accessing file one would need to use </w:t>
      </w:r>
      <w:r>
        <w:t>open() function to let python know what variable we’re going to work with</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">, and what we will be doing with the file. </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">The function returns a “file handle” which is a </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">variable used to perform functions on a file. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Handle=open(filename,mode)</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Fhand=open(</w:t>
      </w:r>
      <w:r>
        <w:t>‘mbox.txt</w:t>
      </w:r>
      <w:r>
        <w:t>’,’r’)</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Handle=returns a handle use to manipulate the file</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Filename= is a string</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Mode= optional and should be ‘r’ if user wants to read file and ‘w’ if user wants to write </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">to the file. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>-</w:t>
      </w:r>
      <w:r>
        <w:t>‘newline” character is used to indicate when a line ends, represented as “</w:t>
      </w:r>
      <w:r>
        <w:t>\n” in strings.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">  A textfile has newlines at the end of each line.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">-a file handle can </w:t>
      </w:r>
      <w:r>
        <w:t>be a sequence of strings and can use the for statement to iterate through a sequence(ordered set)</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:lastRenderedPageBreak/>
        <w:t>Xfile=open(‘mbox.text’)</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>For cheese in xfile:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:tab/>
        <w:t>Print(cheese)</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>-</w:t>
      </w:r>
      <w:r>
        <w:t>We can put an if statement in our for loop to only print lines that meet some c</w:t>
      </w:r>
      <w:r>
        <w:t>riteria.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>-Can skip line by using the ‘continue’ statement</w:t>
      </w:r>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
    </w:p>
</w:body>
'@

$insertionRange = $d.Paragraphs.Item($idx + 1).Range
$insertionRange.InsertXML($newBlockXml)

# Rewrite the anchor paragraph so it no longer carries the _GoBack
# bookmark (which now belongs to the last of the newly inserted
# paragraphs), while keeping its original identity/attributes and runs.
$anchorXml = @'
<w:p w14:paraId="1FD67B21" w14:textId="30B0D507" w:rsidR="00C80C46" w:rsidRDefault="00C80C46"><w:r><w:t>Str.strip(</w:t></w:r><w:r w:rsidR="00F726BF"><w:t>[chars])</w:t></w:r><w:r w:rsidR="00CC3EF6"><w:t xml:space="preserve"> (Both</w:t></w:r><w:r w:rsidR="002D01D6"><w:t xml:space="preserve"> beginning and ending whitespace)</w:t></w:r></w:p>
'@
$d.Paragraphs.Item($idx).Range.InsertXML($anchorXml)

# The original "Str.upper()" / blank / "Search and replace..." paragraphs
# now sit right after the 20 newly inserted paragraphs (duplicated
# content); remove them since that content now lives earlier in the
# document, immediately after the anchor paragraph.
$dupStart = $idx + 1 + 20
$firstDup = $d.Paragraphs.Item($dupStart).Range.Start
$lastDup = $d.Paragraphs.Item($dupStart + 2).Range.End
$dupRange = $d.Range($firstDup, $lastDup)
$dupRange.Delete()

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
